# FinaLytics_Template_2.xlsx - Get_Quotes() kann nun domestic und foreign
# quotes runterladen: update two Yahoo tickers on Sheet1 to their new
# (domestic-exchange) symbols and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 first, then row 6 - keeps the new shared-string order identical to
# the reference edit (MSE.PA appended before DX2J.F).
$ws.Range("A7").Value = "MSE.PA"
$ws.Range("A6").Value = "DX2J.F"

# Selection moved from B12 to E8.
$ws.Range("E8").Select()
